# Organizing QB website 1/7
# - Appends a new "42" reference to the occurrence codes in column E.
# - Rows 2 & 3 get their existing text updated in place (same shared string).
# - Rows 4-7 switch to a new, distinct value "15, 34, 38, 42".
# - The active selection / scroll position moves from B7 to E7 (with C5
#   scrolled to the top-left of the viewport).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = "15, 34, 38. 42"
$ws.Range("E3").Value = "15. 34, 38, 42"
$ws.Range("E4").Value = "15, 34, 38, 42"
$ws.Range("E5").Value = "15, 34, 38, 42"
$ws.Range("E6").Value = "15, 34, 38, 42"
$ws.Range("E7").Value = "15, 34, 38, 42"

# Reposition the view: scroll so C5 is the top-left visible cell and select E7
# (mirrors the sheetView/selection change recorded in the workbook).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("E7").Select()
